# Update crypto price list: refresh Price (D) values and Hora (G) hour stamp
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") changed for rows with a numeric quote; rows holding "--" are left untouched
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "268.86"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "22.84"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "6.369"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06216"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "3.630"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "6.704"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.392"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.8319"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.01369"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1611"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08250"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03403"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.03191"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.09304"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.941"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.001708"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.04859"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.006246"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.005380"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.001093"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0001503"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.757"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.371"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.3342"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0002696"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04661"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006886"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1157"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.003472"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01235"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00006286"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000754"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.7036"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1399"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002111"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.01247"

# Column G ("Hora") advances from 19 to 20 for every data row (2-51)
$ws.Range("G2:G51").NumberFormat = "@"
$ws.Range("G2:G51").Value = "20"
